# Auto-generated Excel COM-interop script to apply the "scheduled runner" data refresh
# to the Ultros_Profits workbook (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC!row9
$ws_ALC.Range("H9").Value = 118.71429
$ws_ALC.Range("I9").Value = 126
$ws_ALC.Range("K9").Value = 126
$ws_ALC.Range("M9").Value = 43

# ALC!row61
$ws_ALC.Range("H61").Value = 568.75
$ws_ALC.Range("I61").Value = 568.75
$ws_ALC.Range("K61").Value = 1706.25
$ws_ALC.Range("M61").Value = -1534.25

# ALC!row70
$ws_ALC.Range("H70").Value = 2616.3242
$ws_ALC.Range("I70").Value = 2135.2222
$ws_ALC.Range("J70").Value = 3072.1052
$ws_ALC.Range("K70").Value = 6405.6666
$ws_ALC.Range("L70").Value = 9216.3156
$ws_ALC.Range("M70").Value = -6135.6666
$ws_ALC.Range("N70").Value = -9756.3156

# ALC!row73
$ws_ALC.Range("H73").Value = 2616.3242
$ws_ALC.Range("I73").Value = 2135.2222
$ws_ALC.Range("J73").Value = 3072.1052
$ws_ALC.Range("K73").Value = 6405.6666
$ws_ALC.Range("L73").Value = 9216.3156
$ws_ALC.Range("M73").Value = -5469.6666
$ws_ALC.Range("N73").Value = -11088.3156

# ALC!row74
$ws_ALC.Range("H74").Value = 7017.4375
$ws_ALC.Range("I74").Value = 5189.9165
$ws_ALC.Range("K74").Value = 5189.9165
$ws_ALC.Range("M74").Value = -4253.9165

# ALC!row77
$ws_ALC.Range("H77").Value = 7017.4375
$ws_ALC.Range("I77").Value = 5189.9165
$ws_ALC.Range("K77").Value = 25949.5825
$ws_ALC.Range("M77").Value = -21269.5825

# ALC!row107
$ws_ALC.Range("H107").Value = 621.25
$ws_ALC.Range("I107").Value = 745
$ws_ALC.Range("J107").Value = 497.5
$ws_ALC.Range("K107").Value = 745
$ws_ALC.Range("L107").Value = 497.5
$ws_ALC.Range("M107").Value = 1175
$ws_ALC.Range("N107").Value = -4337.5

# ALC!row132
$ws_ALC.Range("H132").Value = 13401.439
$ws_ALC.Range("I132").Value = 1226.3235
$ws_ALC.Range("K132").Value = 3678.9705
$ws_ALC.Range("M132").Value = -1148.9705

# ALC!row138
$ws_ALC.Range("H138").Value = 2748.2375
$ws_ALC.Range("I138").Value = 1314.8077
$ws_ALC.Range("J138").Value = 3438.4075
$ws_ALC.Range("K138").Value = 3944.4231
$ws_ALC.Range("L138").Value = 10315.2225
$ws_ALC.Range("M138").Value = 1195.5769
$ws_ALC.Range("N138").Value = -20595.2225

# ARM!row132
$ws_ARM.Range("H132").Value = 3648.3333
$ws_ARM.Range("I132").Value = 3431.0688
$ws_ARM.Range("K132").Value = 10293.2064
$ws_ARM.Range("M132").Value = -7763.206399999999

# BSM!row86
$ws_BSM.Range("H86").Value = 3987.5293
$ws_BSM.Range("I86").Value = 1578.875
$ws_BSM.Range("K86").Value = 1578.875
$ws_BSM.Range("M86").Value = -455.875

# BSM!row89
$ws_BSM.Range("H89").Value = 3987.5293
$ws_BSM.Range("I89").Value = 1578.875
$ws_BSM.Range("K89").Value = 7894.375
$ws_BSM.Range("M89").Value = -2278.375

# BSM!row92
$ws_BSM.Range("H92").Value = 48833.332
$ws_BSM.Range("J92").Value = 48833.332
$ws_BSM.Range("L92").Value = 48833.332
$ws_BSM.Range("N92").Value = -53825.332

# BSM!row94
$ws_BSM.Range("H94").Value = 3962.9285
$ws_BSM.Range("I94").Value = 3046.5454
$ws_BSM.Range("K94").Value = 3046.5454
$ws_BSM.Range("M94").Value = -2595.5454

# BSM!row134
$ws_BSM.Range("H134").Value = 2567.8865
$ws_BSM.Range("I134").Value = 1290.6666
$ws_BSM.Range("K134").Value = 3871.9998
$ws_BSM.Range("M134").Value = -1336.9998

# CRP!row6
$ws_CRP.Range("H6").Value = 2549.5
$ws_CRP.Range("J6").Value = 2549.5
$ws_CRP.Range("L6").Value = 2549.5
$ws_CRP.Range("N6").Value = -2775.5

# CRP!row86
$ws_CRP.Range("H86").Value = 66959.60000000001
$ws_CRP.Range("I86").Value = 104933
$ws_CRP.Range("J86").Value = 9999.5
$ws_CRP.Range("K86").Value = 104933
$ws_CRP.Range("L86").Value = 9999.5
$ws_CRP.Range("M86").Value = -103810
$ws_CRP.Range("N86").Value = -12245.5

# CRP!row89
$ws_CRP.Range("H89").Value = 66959.60000000001
$ws_CRP.Range("I89").Value = 104933
$ws_CRP.Range("J89").Value = 9999.5
$ws_CRP.Range("K89").Value = 524665
$ws_CRP.Range("L89").Value = 49997.5
$ws_CRP.Range("M89").Value = -519049
$ws_CRP.Range("N89").Value = -61229.5

# CRP!row134
$ws_CRP.Range("H134").Value = 2571.1333
$ws_CRP.Range("I134").Value = 1417.4857
$ws_CRP.Range("K134").Value = 4252.4571
$ws_CRP.Range("M134").Value = -1717.4571

# CUL!row4
$ws_CUL.Range("H4").Value = 46991296
$ws_CUL.Range("I4").Value = 59429924
$ws_CUL.Range("J4").Value = 932.2222
$ws_CUL.Range("K4").Value = 178289772
$ws_CUL.Range("L4").Value = 2796.6666
$ws_CUL.Range("M4").Value = -178289660
$ws_CUL.Range("N4").Value = -3020.6666

# CUL!row7
$ws_CUL.Range("H7").Value = 1219.0625
$ws_CUL.Range("I7").Value = 1328.7858
$ws_CUL.Range("J7").Value = 451
$ws_CUL.Range("K7").Value = 3986.3574
$ws_CUL.Range("L7").Value = 1353
$ws_CUL.Range("M7").Value = -3874.3574
$ws_CUL.Range("N7").Value = -1577

# CUL!row17
$ws_CUL.Range("H17").Value = 5048.5
$ws_CUL.Range("I17").Value = 5048.5
$ws_CUL.Range("K17").Value = 15145.5
$ws_CUL.Range("M17").Value = -14976.5

# CUL!row22
$ws_CUL.Range("H22").Value = 37044028
$ws_CUL.Range("I22").Value = 579.6
$ws_CUL.Range("K22").Value = 1738.8
$ws_CUL.Range("M22").Value = -1569.8

# CUL!row27
$ws_CUL.Range("H27").Value = 37044028
$ws_CUL.Range("I27").Value = 579.6
$ws_CUL.Range("K27").Value = 1738.8
$ws_CUL.Range("M27").Value = -1636.8

# CUL!row29
$ws_CUL.Range("H29").Value = 37374468
$ws_CUL.Range("I29").Value = 9723174
$ws_CUL.Range("J29").Value = 111111250
$ws_CUL.Range("K29").Value = 29169522
$ws_CUL.Range("L29").Value = 333333750
$ws_CUL.Range("M29").Value = -29169245
$ws_CUL.Range("N29").Value = -333334304

# CUL!row34
$ws_CUL.Range("H34").Value = 1366.5
$ws_CUL.Range("I34").Value = 479.5
$ws_CUL.Range("J34").Value = 1662.1666
$ws_CUL.Range("K34").Value = 1438.5
$ws_CUL.Range("L34").Value = 4986.4998
$ws_CUL.Range("M34").Value = -1354.5
$ws_CUL.Range("N34").Value = -5154.4998

# CUL!row74
$ws_CUL.Range("H74").Value = 7000
$ws_CUL.Range("J74").Value = 7000
$ws_CUL.Range("L74").Value = 21000
$ws_CUL.Range("N74").Value = -23122

# CUL!row77
$ws_CUL.Range("H77").Value = 7000
$ws_CUL.Range("J77").Value = 7000
$ws_CUL.Range("L77").Value = 63000
$ws_CUL.Range("N77").Value = -73608

# GSM!row18
$ws_GSM.Range("H18").Value = 4999.5
$ws_GSM.Range("I18").Value = 4999
$ws_GSM.Range("K18").Value = 4999
$ws_GSM.Range("M18").Value = -4706

# GSM!row43
$ws_GSM.Range("H43").Value = 7966.6665
$ws_GSM.Range("I43").Value = 933.3333
$ws_GSM.Range("K43").Value = 933.3333
$ws_GSM.Range("M43").Value = -782.3333

# GSM!row46
$ws_GSM.Range("H46").Value = 26875
$ws_GSM.Range("J46").Value = 34833.332
$ws_GSM.Range("L46").Value = 34833.332
$ws_GSM.Range("N46").Value = -35145.332

# GSM!row80
$ws_GSM.Range("H80").Value = 82806.78999999999
$ws_GSM.Range("I80").Value = 161771.42
$ws_GSM.Range("J80").Value = 3842.1428
$ws_GSM.Range("K80").Value = 161771.42
$ws_GSM.Range("L80").Value = 3842.1428
$ws_GSM.Range("M80").Value = -160773.42
$ws_GSM.Range("N80").Value = -5838.1428

# GSM!row83
$ws_GSM.Range("H83").Value = 82806.78999999999
$ws_GSM.Range("I83").Value = 161771.42
$ws_GSM.Range("J83").Value = 3842.1428
$ws_GSM.Range("K83").Value = 808857.1000000001
$ws_GSM.Range("L83").Value = 19210.714
$ws_GSM.Range("M83").Value = -803865.1000000001
$ws_GSM.Range("N83").Value = -29194.714

# GSM!row122
$ws_GSM.Range("H122").Value = 10998.5
$ws_GSM.Range("I122").Value = 10998.5
$ws_GSM.Range("J122").Value = 0
$ws_GSM.Range("K122").Value = 32995.5
$ws_GSM.Range("L122").Value = 0
$ws_GSM.Range("M122").Value = -30545.5
$ws_GSM.Range("N122").ClearContents()

# LTW!row16
$ws_LTW.Range("H16").Value = 125.85714
$ws_LTW.Range("I16").Value = 130.16667
$ws_LTW.Range("J16").Value = 100
$ws_LTW.Range("K16").Value = 130.16667
$ws_LTW.Range("L16").Value = 100
$ws_LTW.Range("M16").Value = 39.83332999999999
$ws_LTW.Range("N16").Value = -440

# LTW!row55
$ws_LTW.Range("H55").Value = 3693.2856
$ws_LTW.Range("I55").Value = 10125
$ws_LTW.Range("K55").Value = 10125
$ws_LTW.Range("M55").Value = -9952

# LTW!row61
$ws_LTW.Range("H61").Value = 2226.7334
$ws_LTW.Range("I61").Value = 783.4583
$ws_LTW.Range("J61").Value = 7999.8335
$ws_LTW.Range("K61").Value = 783.4583
$ws_LTW.Range("L61").Value = 7999.8335
$ws_LTW.Range("M61").Value = -581.4583
$ws_LTW.Range("N61").Value = -8403.833500000001

# LTW!row93
$ws_LTW.Range("H93").Value = 837989.75
$ws_LTW.Range("I93").Value = 6313.1665
$ws_LTW.Range("K93").Value = 6313.1665
$ws_LTW.Range("M93").Value = -5065.1665

# LTW!row113
$ws_LTW.Range("H113").Value = 2226.7334
$ws_LTW.Range("I113").Value = 783.4583
$ws_LTW.Range("J113").Value = 7999.8335
$ws_LTW.Range("K113").Value = 783.4583
$ws_LTW.Range("L113").Value = 7999.8335
$ws_LTW.Range("M113").Value = 1386.5417
$ws_LTW.Range("N113").Value = -12339.8335

# LTW!row132
$ws_LTW.Range("H132").Value = 3565.4
$ws_LTW.Range("I132").Value = 2960.3823
$ws_LTW.Range("J132").Value = 5435.4546
$ws_LTW.Range("K132").Value = 8881.1469
$ws_LTW.Range("L132").Value = 16306.3638
$ws_LTW.Range("M132").Value = -6351.1469
$ws_LTW.Range("N132").Value = -21366.3638

# WVR!row107
$ws_WVR.Range("H107").Value = 583
$ws_WVR.Range("I107").Value = 527.3333
$ws_WVR.Range("K107").Value = 1581.9999
$ws_WVR.Range("M107").Value = 338.0001

# WVR!row111
$ws_WVR.Range("H111").Value = 69998.5
$ws_WVR.Range("J111").Value = 69998.5
$ws_WVR.Range("L111").Value = 69998.5
$ws_WVR.Range("N111").Value = -78178.5

# WVR!row132
$ws_WVR.Range("H132").Value = 2238.614
$ws_WVR.Range("I132").Value = 2026.5385
$ws_WVR.Range("K132").Value = 6079.6155
$ws_WVR.Range("M132").Value = -3549.6155

# WVR!row136
$ws_WVR.Range("H136").Value = 2182.805
$ws_WVR.Range("I136").Value = 1082.9354
$ws_WVR.Range("J136").Value = 5592.4
$ws_WVR.Range("K136").Value = 3248.8062
$ws_WVR.Range("L136").Value = 16777.2
$ws_WVR.Range("M136").Value = -698.8062
$ws_WVR.Range("N136").Value = -21877.2

